$d = $word.ActiveDocument

# --- Locate the anchor text without mutating anything yet -----------------
# "java, salvataggio d..." is a single run in the original document; we need
# to insert ", java Swing" right after "java" and move the (hidden) _GoBack
# bookmark to sit right after that insertion (before ", salvataggio").
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("java, salvataggio", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'java, salvataggio' anchor text"
}
$javaEnd = $findRange.Start + 4

# --- Insert the new text ----------------------------------------------
# (This single text mutation is the only one we perform; Word's run
# coalescing after a text edit merges same-formatted runs in the
# paragraph, so we rebuild the desired run boundaries afterwards using
# bookmark carve-outs, which do not trigger that coalescing.)
$insPoint = $d.Range($javaEnd, $javaEnd)
$newText = ", java Swing"
$insPoint.InsertAfter($newText)

# --- Recompute the split points in the now-longer paragraph ---------------
$spaceBeforeJava = $javaEnd - 4          # just before "java"
$afterNewText    = $javaEnd + $newText.Length   # right after ", java Swing"
$afterOldSplit   = $afterNewText + ", salvataggio d".Length  # before "ati in file a parte"

# --- Re-create the run boundaries that the coalescing merge collapsed -----
$d.Bookmarks.Add("ZzSplit1", $d.Range($spaceBeforeJava, $spaceBeforeJava))
$d.Bookmarks.Add("ZzSplit2", $d.Range($javaEnd, $javaEnd))
$d.Bookmarks.Add("ZzSplit3", $d.Range($afterNewText, $afterNewText))
$d.Bookmarks.Add("ZzSplit4", $d.Range($afterOldSplit, $afterOldSplit))

# --- Move "_GoBack" to sit right after the inserted ", java Swing" --------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Bookmarks("ZzSplit3").Range)

# --- Drop the temporary helper bookmarks (their job -- carving run
#     boundaries -- is already baked into the paragraph's run list) -------
$d.Bookmarks("ZzSplit1").Delete()
$d.Bookmarks("ZzSplit2").Delete()
$d.Bookmarks("ZzSplit3").Delete()
$d.Bookmarks("ZzSplit4").Delete()
